$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57:66 down to 58:67
$ws.Rows.Item(57).Insert()

# Fill in the new row 57 with the new weekly record
$ws.Cells.Item(57, 1).Value = 4
$ws.Cells.Item(57, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value = "Los Lagos"
$ws.Cells.Item(57, 4).Value = 44995
$ws.Cells.Item(57, 5).Value = 10
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100101
$ws.Cells.Item(57, 8).Value = "Berries"
$ws.Cells.Item(57, 9).Value = 100101001
$ws.Cells.Item(57, 10).Value = "Arándano (blue)"
$ws.Cells.Item(57, 11).Value = "Sin especificar"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 200
$ws.Cells.Item(57, 14).Value = 2300
$ws.Cells.Item(57, 15).Value = 2500
$ws.Cells.Item(57, 16).Value = 2400
$ws.Cells.Item(57, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(57, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(57, 19).Value = 1200
$ws.Cells.Item(57, 20).Value = 2
